$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.135.87"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "3.548.26"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.61"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.20%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.214"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.652"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.81"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("E12").Value = "  -4.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.53"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").Value = "4.118.80"
$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "632.28"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +10.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.96"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.70%  "

$ws.Range("D17").Value = "70.180.58"
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.92"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("D19").Value = "3.533.97"
$ws.Range("E19").Value = "  -1.98%  "

$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.56"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.63%  "

$ws.Range("E23").Value = "  +1.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.27"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.94"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.03"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.74"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.31"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.71"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.69"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +18.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.26"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "529.64"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.402"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.25"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.99%  "

$ws.Range("D40").Value = "0.0₃0780"
$ws.Range("E40").Value = "  -1.65%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.53"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.70%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.504.65"
$ws.Range("E42").Value = "  +3.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.137"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0459"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.05%  "

$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("E46").Value = "  +4.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.36"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("E50").Value = "  -3.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.09%  "

